$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing glossary table (Table1), which currently spans A1:D46.
$lo = $ws.ListObjects.Item(1)

# New glossary rows to append: Topic / Term / SeeAlso / Definition
$newRows = @(
    @{ Topic = 'Power';  Term = 'Nema L2-20'; SeeAlso = ''; Definition = 'This is a heavy duty 20 Amp  connector with circular locking contacts  (hence the **L** in the name). There is an *R* and a *P* suffix for the *recepacle* and *plug*.' },
    @{ Topic = 'Power '; Term = 'Nema 5-15';   SeeAlso = ''; Definition = 'This is a regular three-prong   connector  good for 15 Amps . There is an *R* and a *P* suffix for the *recepacle* and *plug*.' },
    @{ Topic = 'Power '; Term = 'Nema 6-20';   SeeAlso = ''; Definition = 'This looks like a regular three-prong plug but one of the blades is turned 90 degrees. This connector is good for 20 Amp.  There is an *R* and a *P* suffix for the *recepacle* and *plug*.' },
    @{ Topic = 'Power';  Term = 'Cam-Lock';    SeeAlso = ''; Definition = 'A power connector type used for large current power distribution (> 100 Amp). See https://en.wikipedia.org/wiki/Camlock_(electrical)' }
)

foreach ($row in $newRows) {
    $newListRow = $lo.ListRows.Add()
    $newRange = $newListRow.Range
    $rowNum = $newRange.Row

    # Copy formatting (styles, wrap, vertical alignment) from the last existing data row.
    $ws.Range("A$($rowNum - 1):D$($rowNum - 1)").Copy()
    $ws.Range("A$($rowNum):D$($rowNum)").PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $row.Topic
    $ws.Cells.Item($rowNum, 2).Value = $row.Term
    if ($row.SeeAlso -ne '') {
        $ws.Cells.Item($rowNum, 3).Value = $row.SeeAlso
    }
    $ws.Cells.Item($rowNum, 4).Value = $row.Definition

    $ws.Rows.Item($rowNum).RowHeight = 34
}

$excel.CutCopyMode = $false

# Match the author's final viewport/selection state.
[void]$ws.Range("D40").Select()
